# PIM-6472: Fix wrong family variant codes
#
# The "code" column (A) of the family_variant sheet contained codes that
# collided with / duplicated other fixtures. Rename them so they are
# unique ("another_..." prefix), matching the corrected fixture data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "another_clothing_color_and_size"
$ws.Range("A3").Value = "another_shoes_size"
$ws.Range("A4").Value = "another_clothing_color_size"

# Tidy up: C4 carried a redundant/duplicate cell style (applyFont=true but
# same font as the default) - reset it back to the sheet's default style,
# same as the rest of the column, by pasting the formatting already used
# elsewhere in the sheet (keeps the cell's value/content untouched).
$ws.Range("A4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Selection moved to cover the whole used range with the active cell back
# at A1 (was sitting on a stale E20 selection outside the used range).
$ws.Range("A1:I4").Select() | Out-Null

# Cosmetic: tab-bar split ratio tweak from the original edit.
$excel.ActiveWindow.TabRatio = 0.5
